# Fruta / hortaliza, semanal
# Insert 3 new observation rows for Cereza (Lapins/Lapins/Santina) just above
# the existing row 206, shifting the rest of the table (old rows 206-258)
# down to rows 209-261, then populate the 3 newly-inserted rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three blank rows at 206-208; everything currently at row 206 and
# below (through 258) shifts down to 209-261, carrying its formatting along.
$ws.Range("A206:A208").EntireRow.Insert()

# --- New row 206 ---
$ws.Range("A206").Value = 9
$ws.Range("B206").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C206").Value = "Metropolitana"
$ws.Range("D206").Value = 44559
$ws.Range("E206").Value = 13
$ws.Range("F206").Value = "Fruta"
$ws.Range("G206").Value = 100103
$ws.Range("H206").Value = "Frutos de hueso (carozo)"
$ws.Range("I206").Value = 100103001
$ws.Range("J206").Value = "Cereza"
$ws.Range("K206").Value = "Lapins"
$ws.Range("L206").Value = "Primera"
$ws.Range("M206").Value = 440
$ws.Range("N206").Value = 9000
$ws.Range("O206").Value = 9000
$ws.Range("P206").Value = 9000
$ws.Range("Q206").Value = "`$/caja 18 kilos"
$ws.Range("R206").Value = "Región de O'Higgins"
$ws.Range("S206").Value = 500
$ws.Range("T206").Value = 18

# --- New row 207 ---
$ws.Range("A207").Value = 9
$ws.Range("B207").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C207").Value = "Metropolitana"
$ws.Range("D207").Value = 44559
$ws.Range("E207").Value = 13
$ws.Range("F207").Value = "Fruta"
$ws.Range("G207").Value = 100103
$ws.Range("H207").Value = "Frutos de hueso (carozo)"
$ws.Range("I207").Value = 100103001
$ws.Range("J207").Value = "Cereza"
$ws.Range("K207").Value = "Lapins"
$ws.Range("L207").Value = "Segunda"
$ws.Range("M207").Value = 380
$ws.Range("N207").Value = 7200
$ws.Range("O207").Value = 7200
$ws.Range("P207").Value = 7200
$ws.Range("Q207").Value = "`$/caja 18 kilos"
$ws.Range("R207").Value = "Región de O'Higgins"
$ws.Range("S207").Value = 400
$ws.Range("T207").Value = 18

# --- New row 208 ---
$ws.Range("A208").Value = 9
$ws.Range("B208").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C208").Value = "Metropolitana"
$ws.Range("D208").Value = 44559
$ws.Range("E208").Value = 13
$ws.Range("F208").Value = "Fruta"
$ws.Range("G208").Value = 100103
$ws.Range("H208").Value = "Frutos de hueso (carozo)"
$ws.Range("I208").Value = 100103001
$ws.Range("J208").Value = "Cereza"
$ws.Range("K208").Value = "Santina"
$ws.Range("L208").Value = "Primera"
$ws.Range("M208").Value = 790
$ws.Range("N208").Value = 4000
$ws.Range("O208").Value = 4500
$ws.Range("P208").Value = 4241
$ws.Range("Q208").Value = "`$/caja 10 kilos"
$ws.Range("R208").Value = "Provincia de Curicó"
$ws.Range("S208").Value = 424
$ws.Range("T208").Value = 10
